$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.899.64"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").Value = "3.477.10"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'601.02"
$ws.Range("E5").Value = "  -3.17%  "
$ws.Range("D6").Value = "'147.75"
$ws.Range("E6").Value = "  -4.74%  "
$ws.Range("D7").Value = "3.473.49"
$ws.Range("E7").Value = "  -2.59%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("D11").Value = "'7.62"
$ws.Range("E11").Value = "  +3.15%  "
$ws.Range("D12").Value = "'0.422"
$ws.Range("E12").Value = "  -3.73%  "
$ws.Range("D13").Value = "'0.0000212"
$ws.Range("E13").Value = "  -4.47%  "
$ws.Range("D14").Value = "4.065.86"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").Value = "'31.21"
$ws.Range("E15").Value = "  -6.13%  "
$ws.Range("D16").Value = "3.471.04"
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").Value = "66.948.28"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  -5.19%  "
$ws.Range("D20").Value = "'15.25"
$ws.Range("E20").Value = "  -4.76%  "
$ws.Range("D21").Value = "'10.04"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'433.25"
$ws.Range("E22").Value = "  -4.95%  "
$ws.Range("D23").Value = "'0.604"
$ws.Range("E23").Value = "  -6.03%  "
$ws.Range("D24").Value = "'78.96"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "3.616.06"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").Value = "'0.0000120"
$ws.Range("E27").Value = "  -8.55%  "
$ws.Range("D28").Value = "'9.78"
$ws.Range("E28").Value = "  -7.33%  "
$ws.Range("D29").Value = "'8.37"
$ws.Range("E29").Value = "  -8.98%  "
$ws.Range("D30").Value = "'2.47"
$ws.Range("E30").Value = "  -3.78%  "
$ws.Range("D31").Value = "'1.58"
$ws.Range("E31").Value = "  -7.73%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "'25.25"
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D35").Value = "3.469.09"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").Value = "'5.90"
$ws.Range("E36").Value = "  -7.27%  "
$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = "  -6.43%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "'7.87"
$ws.Range("E39").Value = "  -4.59%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'173.61"
$ws.Range("E41").Value = "  -4.34%  "
$ws.Range("D42").Value = "'0.0881"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("D43").Value = "'2.08"
$ws.Range("E43").Value = "  -13.15%  "
$ws.Range("D44").Value = "'5.38"
$ws.Range("E44").Value = "  -3.78%  "
$ws.Range("D45").Value = "'0.893"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "'46.42"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "'28.85"
$ws.Range("E47").Value = "  -7.49%  "
$ws.Range("E48").Value = "  -7.27%  "
$ws.Range("D49").Value = "'7.43"
$ws.Range("E49").Value = "  -4.62%  "
$ws.Range("D50").Value = "'2.40"
$ws.Range("E50").Value = "  -10.22%  "
$ws.Range("D51").Value = "'0.971"
$ws.Range("E51").Value = "  -4.80%  "
